{"js": "// Apply the documented text replacements to the course-template document:\n//  1. \"Cours de la formation g\u00e9n\u00e9rale 111.GE\" -> \"111.GE Cours de la formation g\u00e9n\u00e9rale\"\n//  2. \"Philosophie et rationalit\u00e9\"             -> \"Titre du cours :Philosophie et rationalit\u00e9\"\n//  3. \"340-101-MQ\"                             -> \"Numero du cours :340-101-MQ\"\n//  4. \"3-1-3\"                                  -> \"Pond\u00e9ration :3-1-3\"\n//  5. \"1.000000\"                               -> \"Nombre d'unit\u00e9(s) :1\"\n\nconst replacements = [\n  [\"Cours de la formation g\u00e9n\u00e9rale 111.GE\", \"111.GE Cours de la formation g\u00e9n\u00e9rale\"],\n  [\"Philosophie et rationalit\u00e9\", \"Titre du cours :Philosophie et rationalit\u00e9\"],\n  [\"340-101-MQ\", \"Numero du cours :340-101-MQ\"],\n  [\"3-1-3\", \"Pond\u00e9ration :3-1-3\"],\n  [\"1.000000\", \"Nombre d'unit\u00e9(s) :1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the documented text replacements to the course-template document:\n#  1. \"Cours de la formation g\u00e9n\u00e9rale 111.GE\" -> \"111.GE Cours de la formation g\u00e9n\u00e9rale\"\n#  2. \"Philosophie et rationalit\u00e9\"             -> \"Titre du cours :Philosophie et rationalit\u00e9\"\n#  3. \"340-101-MQ\"                             -> \"Numero du cours :340-101-MQ\"\n#  4. \"3-1-3\"                                  -> \"Pond\u00e9ration :3-1-3\"\n#  5. \"1.000000\"                               -> \"Nombre d'unit\u00e9(s) :1\"\n#\n# NOTE: we replace via the *whole paragraph* Range.Text (the Paragraphs\n# collection walks body paragraphs AND table-cell paragraphs) instead of a\n# Find/Replace on a sub-run Range; this keeps each run's existing\n# formatting (empty <w:rPr/>, xml:space=\"preserve\") and avoids Word's\n# smart-quote autocorrect mangling the apostrophe in \"d'unit\u00e9(s)\".\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"Cours de la formation g\u00e9n\u00e9rale 111.GE\"; new = \"111.GE Cours de la formation g\u00e9n\u00e9rale\" },\n    @{ old = \"Philosophie et rationalit\u00e9\";             new = \"Titre du cours :Philosophie et rationalit\u00e9\" },\n    @{ old = \"340-101-MQ\";                             new = \"Numero du cours :340-101-MQ\" },\n    @{ old = \"3-1-3\";                                  new = \"Pond\u00e9ration :3-1-3\" },\n    @{ old = \"1.000000\";                               new = \"Nombre d'unit\u00e9(s) :1\" }\n)\n\nforeach ($rep in $replacements) {\n    foreach ($p in $d.Paragraphs) {\n        # Paragraphs(i).Range.Text includes the trailing paragraph-mark\n        # (chr 13) and, for a table-cell paragraph, the cell-mark (chr 7)\n        # too - strip both before comparing against the plain target text.\n        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($ptext -eq $rep.old) {\n            $p.Range.Text = $rep.new\n            break\n        }\n    }\n}\n"}
